# Sprint 0 presentation update:
# Insert three new content slides into the deck while leaving the
# existing slides' own content untouched (only their position shifts).
#
#   pos2 (new): "Wat is het Probleem?"
#   pos3 (new): "Mogelijke oplossingen?"
#   pos9 (new): "Samenwerkingsovereenkomst"   (inserted right after
#               "Wat hebben we gedaan?" and before "Volgende Sprint")

$p = $ppt.ActivePresentation

# Use the existing "Vragen?" slide (Title + Content placeholder layout,
# nl-NL formatted) as the template for new slides so the new slides come
# out with the same placeholder/formatting shape as the rest of the deck.
$template = $p.Slides.Item($p.Slides.Count)

# --- New slide: "Wat is het Probleem?" -> becomes slide 2 -----------------
$range1 = $template.Duplicate()
$newSlide1 = $range1.Item(1)
$newSlide1.MoveTo(2)
$newSlide1.Shapes.Item(1).TextFrame.TextRange.Text = "Wat is het Probleem?"
$newSlide1.Shapes.Item(2).TextFrame.TextRange.Text = "Floating farm`nBalans`nInstabiliteit door Koeien"

# --- New slide: "Mogelijke oplossingen?" -> becomes slide 3 ---------------
$range2 = $template.Duplicate()
$newSlide2 = $range2.Item(1)
$newSlide2.MoveTo(3)
$newSlide2.Shapes.Item(1).TextFrame.TextRange.Text = "Mogelijke oplossingen?"
$newSlide2.Shapes.Item(2).TextFrame.TextRange.Text = "Ballasttanks (meerdere varianten)`nNiet zeker want…"

# --- New slide: "Samenwerkingsovereenkomst" -> becomes slide 9 ------------
# (after the two insertions above, the original slides shifted down by two,
#  so "Wat hebben we gedaan?" is now slide 8 and "Volgende Sprint" slide 9)
$range3 = $template.Duplicate()
$newSlide3 = $range3.Item(1)
$newSlide3.MoveTo(9)
$newSlide3.Shapes.Item(1).TextFrame.TextRange.Text = "Samenwerkingsovereenkomst"
$newSlide3.Shapes.Item(2).TextFrame.TextRange.Text = "Afspraken`nSterke- en Zwakke punten`nRol/Taakverdeling`n"

Write-Output ("Final slide count: " + $p.Slides.Count)
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    Write-Output ($i.ToString() + ": " + $s.Shapes.Item(1).TextFrame.TextRange.Text)
}
